# Weekly update of the "Betarraga" price series: the newest week's figures
# get inserted at the top of the data block (rows 168:169), every
# already-recorded week shifts down two rows, and the oldest week that falls
# off the bottom is re-appended as new rows 190:191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-append the current last week (rows 188:189) as the new rows 190:191
#    before anything below gets overwritten by the shift in step 2.
$ws.Range("A188:R189").Copy($ws.Range("A190:R191"))

# 2) Shift the existing weekly blocks (rows 168:187) down by one block
#    (two rows) into 170:189, pushing older data further down the sheet.
$ws.Range("A168:R187").Copy($ws.Range("A170:R189"))

# 3) Write this week's new reporting date into the now-vacated top block.
$ws.Range("D168:D169").Value = 44449
